$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "291.47"
Set-TextValue $ws.Range("E2") "0.34%"
Set-TextValue $ws.Range("D3") "31.04"
Set-TextValue $ws.Range("E3") "0.85%"
Set-TextValue $ws.Range("D4") "4.950"
Set-TextValue $ws.Range("E4") "1.71%"
Set-TextValue $ws.Range("E5") "2.65%"
Set-TextValue $ws.Range("D6") "2.233"
Set-TextValue $ws.Range("E6") "-8.16%"
Set-TextValue $ws.Range("D7") "7.720"
Set-TextValue $ws.Range("E7") "0.99%"
Set-TextValue $ws.Range("D8") "0.9183"
Set-TextValue $ws.Range("E8") "2.47%"
Set-TextValue $ws.Range("D9") "0.09408"
Set-TextValue $ws.Range("E9") "16.74%"
Set-TextValue $ws.Range("D10") "0.1722"
Set-TextValue $ws.Range("E10") "3.13%"
Set-TextValue $ws.Range("D11") "0.08329"
Set-TextValue $ws.Range("E11") "2.15%"
Set-TextValue $ws.Range("D12") "0.03220"
Set-TextValue $ws.Range("E12") "4.81%"
Set-TextValue $ws.Range("D13") "0.09978"
Set-TextValue $ws.Range("E13") "-0.53%"
Set-TextValue $ws.Range("D14") "0.001493"
Set-TextValue $ws.Range("E14") "-0.29%"
Set-TextValue $ws.Range("B15") "TigerCash"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.005756"
Set-TextValue $ws.Range("E15") "0.30%"
Set-TextValue $ws.Range("B16") "LEO"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.475"
Set-TextValue $ws.Range("E16") "-0.13%"
Set-TextValue $ws.Range("B17") "GateToken"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "3.746"
Set-TextValue $ws.Range("E17") "1.27%"
Set-TextValue $ws.Range("B18") "BTSEToken"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D18") "2.128"
Set-TextValue $ws.Range("E18") "2.53%"
Set-TextValue $ws.Range("B19") "BitpandaEcosystemToken"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D19") "0.3329"
Set-TextValue $ws.Range("E19") "0.47%"
Set-TextValue $ws.Range("B20") "ProBitToken"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D20") "0.1301"
Set-TextValue $ws.Range("E20") "1.07%"
Set-TextValue $ws.Range("B21") "MCDex"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D21") "4.169"
Set-TextValue $ws.Range("E21") "5.03%"
Set-TextValue $ws.Range("B22") "ZBToken"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws.Range("D22") "0.2115"
Set-TextValue $ws.Range("E22") "0.30%"
Set-TextValue $ws.Range("B23") "CoinExToken"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D23") "0.04530"
Set-TextValue $ws.Range("E23") "0.27%"
Set-TextValue $ws.Range("D24") "0.001214"
Set-TextValue $ws.Range("E24") "0.13%"
Set-TextValue $ws.Range("D25") "0.004255"
Set-TextValue $ws.Range("E25") "-3.43%"
Set-TextValue $ws.Range("E26") "-0.49%"
Set-TextValue $ws.Range("D27") "0.0003380"
Set-TextValue $ws.Range("E27") "-0.49%"
Set-TextValue $ws.Range("D39") "0.01606"
Set-TextValue $ws.Range("E39") "1.13%"
Set-TextValue $ws.Range("E40") "4.30%"
Set-TextValue $ws.Range("D41") "0.007427"
Set-TextValue $ws.Range("E41") "2.43%"
Set-TextValue $ws.Range("D42") "0.009813"
Set-TextValue $ws.Range("E42") "-2.00%"
Set-TextValue $ws.Range("D43") "0.1354"
Set-TextValue $ws.Range("E43") "3.14%"
Set-TextValue $ws.Range("D44") "0.002151"
Set-TextValue $ws.Range("E44") "6.23%"
Set-TextValue $ws.Range("D45") "0.009636"
Set-TextValue $ws.Range("E45") "5.06%"
Set-TextValue $ws.Range("D46") "0.00006101"
Set-TextValue $ws.Range("E46") "6.84%"
Set-TextValue $ws.Range("D47") "0.00000000747"
Set-TextValue $ws.Range("E47") "-0.46%"
Set-TextValue $ws.Range("D48") "2.631"
Set-TextValue $ws.Range("E48") "17.42%"
Set-TextValue $ws.Range("D49") "0.001992"
Set-TextValue $ws.Range("E49") "-31.34%"
Set-TextValue $ws.Range("D50") "0.00002092"
Set-TextValue $ws.Range("E50") "-0.46%"
Set-TextValue $ws.Range("D51") "0.0001992"
Set-TextValue $ws.Range("E51") "-0.46%"
